# ------------------------------------------------------------------
# Applies the "output generated at 456a3b4" update to 北京-漫展信息.xlsx
# Sheets (by tab name):
#   展览     (sheet1) - Exhibitions   - value-only updates to column F
#   演出     (sheet2) - Performances  - value updates + 2 new rows
#   本地生活 (sheet3) - Local life    - value-only updates + 1 new row
#   全部类型 (sheet4) - All types     - value-only updates to column F
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet "展览" (Exhibitions) - update "想去人数" (column F) values
# ============================================================
$wsExpo = $wb.Worksheets.Item("展览")

$expoUpdates = @{
    5  = 8150
    8  = 89
    9  = 7084
    11 = 546
    12 = 489
    14 = 706
    21 = 69
    22 = 11617
    23 = 5
    24 = 129
    25 = 2265
    27 = 3164
    29 = 2687
    32 = 285
    33 = 46
    35 = 1612
    36 = 74
    37 = 105
    38 = 5811
    39 = 80
    40 = 1789
    41 = 1244
    42 = 843
    43 = 159
    46 = 1103
    47 = 1519
    48 = 99
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value2 = $expoUpdates[$row]
}

# ============================================================
# Sheet "演出" (Performances)
# ============================================================
$wsShow = $wb.Worksheets.Item("演出")

# --- plain value fixes (no row shift) ---
$wsShow.Cells.Item(8, 6).Value2 = 250
$wsShow.Cells.Item(20, 6).Value2 = 66

# --- insert a new row 21 for "北京·《山丘》音乐教父 经典情歌金曲翻唱演唱会" ---
# (pushes the old rows 21-25 down to 22-26)
$wsShow.Rows.Item(21).Insert()

# Re-apply the index-column style (bold/centered/bordered) used by every
# other row, since Insert() otherwise invents a slightly different style.
$wsShow.Cells.Item(20, 1).Copy()
$wsShow.Cells.Item(21, 1).PasteSpecial(-4122)  # xlPasteFormats
$wsShow.Application.CutCopyMode = $false

$wsShow.Cells.Item(21, 1).Value2 = 20

# Force the "yyyy-mm-dd" text to stay literal text (not auto-parsed into a
# date serial) the same way the other rows store it, then drop the
# temporary text format so the cell keeps the sheet's default (no) style.
$wsShow.Cells.Item(21, 2).NumberFormat = "@"
$wsShow.Cells.Item(21, 2).Value2 = "2024-08-23"
$wsShow.Cells.Item(21, 2).ClearFormats()

$wsShow.Cells.Item(21, 3).Value2 = "北京·《山丘》音乐教父 经典情歌金曲翻唱演唱会"
$wsShow.Cells.Item(21, 4).Value2 = "大江胡同121号2幢负1层 北京门空间 TheDoorLiveHouse"
$wsShow.Cells.Item(21, 5).Value2 = "2024.08.23 19:30-08.23 21:00"
$wsShow.Cells.Item(21, 6).Value2 = 0
$wsShow.Cells.Item(21, 7).Value2 = 98
$wsShow.Cells.Item(21, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89358"
$wsShow.Cells.Item(21, 9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/noqwx8Qu1721116074567.jpeg"

# --- renumber the index column for the rows that shifted down ---
$wsShow.Cells.Item(22, 1).Value2 = 21
$wsShow.Cells.Item(23, 1).Value2 = 22
$wsShow.Cells.Item(24, 1).Value2 = 23
$wsShow.Cells.Item(25, 1).Value2 = 24
$wsShow.Cells.Item(26, 1).Value2 = 25

# --- append a new row 27 for "北京·伦敦西区音乐剧明星演唱会-经典版" ---
$wsShow.Rows.Item(27).Insert()

$wsShow.Cells.Item(26, 1).Copy()
$wsShow.Cells.Item(27, 1).PasteSpecial(-4122)  # xlPasteFormats
$wsShow.Application.CutCopyMode = $false

$wsShow.Cells.Item(27, 1).Value2 = 26

$wsShow.Cells.Item(27, 2).NumberFormat = "@"
$wsShow.Cells.Item(27, 2).Value2 = "2024-10-25"
$wsShow.Cells.Item(27, 2).ClearFormats()

$wsShow.Cells.Item(27, 3).Value2 = "北京·伦敦西区音乐剧明星演唱会-经典版"
$wsShow.Cells.Item(27, 4).Value2 = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$wsShow.Cells.Item(27, 5).Value2 = "2024.10.25 19:30-10.26 21:30"
$wsShow.Cells.Item(27, 6).Value2 = 0
$wsShow.Cells.Item(27, 7).Value2 = 144
$wsShow.Cells.Item(27, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89359"
$wsShow.Cells.Item(27, 9).Value2 = "//i0.hdslb.com/bfs/openplatform/202407/PzPiEKUI1721114840552.jpeg"

# ============================================================
# Sheet "本地生活" (Local life)
# ============================================================
$wsLocal = $wb.Worksheets.Item("本地生活")

# --- plain value fixes (no row shift) ---
$wsLocal.Cells.Item(2, 6).Value2 = 240
$wsLocal.Cells.Item(3, 6).Value2 = 377

# --- append a new row 4 for "北京·一起幻城动漫游戏嘉年华" ---
$wsLocal.Rows.Item(4).Insert()

$wsLocal.Cells.Item(3, 1).Copy()
$wsLocal.Cells.Item(4, 1).PasteSpecial(-4122)  # xlPasteFormats
$wsLocal.Application.CutCopyMode = $false

$wsLocal.Cells.Item(4, 1).Value2 = 3

$wsLocal.Cells.Item(4, 2).NumberFormat = "@"
$wsLocal.Cells.Item(4, 2).Value2 = "2024-07-27"
$wsLocal.Cells.Item(4, 2).ClearFormats()

$wsLocal.Cells.Item(4, 3).Value2 = "北京·一起幻城动漫游戏嘉年华"
$wsLocal.Cells.Item(4, 4).Value2 = "小关路39号 北投购物公园"
$wsLocal.Cells.Item(4, 5).Value2 = "2024.07.27 10:00-08.08 22:00"
$wsLocal.Cells.Item(4, 6).Value2 = 1
$wsLocal.Cells.Item(4, 7).Value2 = 79
$wsLocal.Cells.Item(4, 8).Value2 = "https://show.bilibili.com/platform/detail.html?id=89323"
$wsLocal.Cells.Item(4, 9).Value2 = "//i1.hdslb.com/bfs/openplatform/202407/uKBmLDLW1721043966929.jpeg"

# ============================================================
# Sheet "全部类型" (All types) - update "想去人数" (column F) values
# ============================================================
$wsAll = $wb.Worksheets.Item("全部类型")

$allUpdates = @{
    3  = 240
    4  = 377
    7  = 8150
    10 = 89
    11 = 7084
    12 = 7084
    14 = 546
    15 = 489
    16 = 706
    21 = 250
    22 = 69
    25 = 11618
    27 = 5
    28 = 129
    29 = 2265
    30 = 2265
    31 = 3164
    32 = 2687
    34 = 285
    35 = 46
    38 = 1612
    39 = 74
    40 = 105
    41 = 5811
    42 = 66
    43 = 1789
    45 = 1244
    46 = 843
    49 = 1103
    50 = 1519
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value2 = $allUpdates[$row]
}

$wb.Save()
